$wb = $excel.ActiveWorkbook

# Rename second sheet from strategy_id-6000 to strategy_id-6002
$ws2 = $wb.Worksheets.Item("strategy_id-6000")
$ws2.Name = "strategy_id-6002"

# Sheet 1: strategy_id-0, clear C6 (remove variable_trajectory_group value)
$ws1 = $wb.Worksheets.Item("strategy_id-0")
$ws1.Range("C6").ClearContents()

# Sheet 2 (now strategy_id-6002): clear C2, and update P2:AR2 with new curve values
$ws2.Range("C2").ClearContents()

$ws2.Range("P2").Value = 0.9253333333333333
$ws2.Range("Q2").Value = 0.9013333333333334
$ws2.Range("R2").Value = 0.8780000000000001
$ws2.Range("S2").Value = 0.8553333333333334
$ws2.Range("T2").Value = 0.8333333333333334
$ws2.Range("U2").Value = 0.8120000000000001
$ws2.Range("V2").Value = 0.7913333333333333
$ws2.Range("W2").Value = 0.7713333333333333
$ws2.Range("X2").Value = 0.752
$ws2.Range("Y2").Value = 0.7333333333333334
$ws2.Range("Z2").Value = 0.7153333333333333
$ws2.Range("AA2").Value = 0.698
$ws2.Range("AB2").Value = 0.6813333333333333
$ws2.Range("AC2").Value = 0.6653333333333333
$ws2.Range("AD2").Value = 0.65
$ws2.Range("AE2").Value = 0.6353333333333333
$ws2.Range("AF2").Value = 0.6213333333333333
$ws2.Range("AG2").Value = 0.6080000000000001
$ws2.Range("AH2").Value = 0.5953333333333333
$ws2.Range("AI2").Value = 0.5833333333333333
$ws2.Range("AJ2").Value = 0.5720000000000001
$ws2.Range("AK2").Value = 0.5613333333333334
$ws2.Range("AL2").Value = 0.5513333333333333
$ws2.Range("AM2").Value = 0.542
$ws2.Range("AN2").Value = 0.5333333333333333
$ws2.Range("AO2").Value = 0.5266666666666666
$ws2.Range("AP2").Value = 0.52
$ws2.Range("AQ2").Value = 0.5133333333333333
$ws2.Range("AR2").Value = 0.5066666666666667
